$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The oldest year row (row 2, "2009年") is removed entirely; the following
# row ("2010年", row 3) shifts up to take its place as row 2.
$ws.Rows.Item(2).Delete()
